$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Row 2 / cell G2: reformat the insert (POST) JSON body onto multiple lines
# ---------------------------------------------------------------------------
$g2 = @'
{
"doorBanSn":"1111",
"ip":"1.1.1.1",
"doorBanName":"门禁1",
"manufacturer":"厂商1",
"model":"型号1",
"longitude":120.333,
"latitude":20.333,
"regionCode":"330104",
"communityCode":"330104",
"direction":"0",
"installationAddress":"杭州下去",
"remark":"没有备注",
"state":0
}
'@
$ws.Range("G2").Value = $g2

# ---------------------------------------------------------------------------
# 2. Row 4 / cell D4: append the "?id=66" query string to the DELETE url
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "http://181.181.0.33:22020/api/acs/v1/door_ban/?id=66"

# ---------------------------------------------------------------------------
# 3. Row 5 / cell G5: replace the PUT (update) JSON body with the new payload
# ---------------------------------------------------------------------------
$g5 = @'
{
"id":"68",
"doorBanSn":"2222",
"ip":"2.2.2.2",
"doorBanName":"门禁1",
"manufacturer":"厂商1",
"model":"型号1",
"longitude":120.333,
"latitude":20.333,
"regionCode":"330104",
"communityCode":"330104",
"direction":"0",
"installationAddress":"杭州下去",
"remark":"没有备注",
"state":0
}
'@
$ws.Range("G5").Value = $g5

# ---------------------------------------------------------------------------
# 4. Row heights changed because of the re-wrapped / replaced text
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 216
$ws.Rows.Item(5).RowHeight = 229.5

# ---------------------------------------------------------------------------
# 5. G5 switches from the "data" style to the plain wrap style (same as G2/G4)
# ---------------------------------------------------------------------------
$ws.Range("G5").Style = $ws.Range("G2").Style

# ---------------------------------------------------------------------------
# 6. Row 4 gains an (empty) G4 cell formatted like the rest of the row
# ---------------------------------------------------------------------------
$ws.Range("G4").Style = $ws.Range("B4").Style

# ---------------------------------------------------------------------------
# 7. Update hyperlink captions for D4 (the underlying target keeps pointing at
#    the same relationship; only the visible text/tooltip change) while
#    leaving D2/D3/D5 untouched in content (re-applied identically so the
#    whole collection stays internally consistent).
# ---------------------------------------------------------------------------
$ws.Range("D4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "http://181.181.0.33:22020/api/acs/v1/door_ban/insert", "", "http://181.181.0.33:22020/api/acs/v1/door_ban/insert", "http://181.181.0.33:22020/api/acs/v1/door_ban/insert") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "http://181.181.0.33:22020/api/acs/v1/door_ban/", "", "http://181.181.0.33:22020/api/acs/v1/door_ban/?id=66", "http://181.181.0.33:22020/api/acs/v1/door_ban/?id=66") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "http://181.181.0.33:22020/api/acs/v1/door_ban/state", "", "http://181.181.0.33:22020/api/acs/v1/door_ban/state", "http://181.181.0.33:22020/api/acs/v1/door_ban/state") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "http://181.181.0.33:22020/api/acs/v1/door_ban/update", "", "http://181.181.0.166:22020/api/acs/v1/door_ban/update", "http://181.181.0.33:22020/api/acs/v1/door_ban/update") | Out-Null

# ---------------------------------------------------------------------------
# 8. Move the active selection to E5 (from E2)
# ---------------------------------------------------------------------------
$ws.Range("E5").Select()
